$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert two new paragraphs at the very beginning of the document:
#    "Scientific background" (bold heading) and the "Lie detection..." body.
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs(1)
$insertPoint = $d.Range($firstPara.Range.Start, $firstPara.Range.Start)

$newParasXml = '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Scientific background</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Lie detection via video leverages involuntary facial expressions, eye movements, and body gestures potentially indicative of deception. 3D Convolutional Neural Networks (3D CNNs) excel in this area by processing spatial and temporal data together, capturing transient expressions and movements critical in identifying lies. They analyze motion and expressions across frames, detecting patterns and inconsistencies. On the other hand, CNN-LSTMs combine CNNs'' spatial feature extraction with LSTMs'' ability to manage long sequences, making them adept at recognizing patterns and contextual behaviors over time. </w:t></w:r></w:p>'

$pkgXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $newParasXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($pkgXml)

# ---------------------------------------------------------------------------
# 2. "outline" -> "outlined"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("outline in the following major steps", $true, $false, $false, $false, $false, $true, 1, $false, "outlined in the following major steps", 2)

# ---------------------------------------------------------------------------
# 3. Remove "(and verify that the extracted image is a face)" parenthetical
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(". (and verify that the extracted image is a face)", $true, $false, $false, $false, $false, $true, 1, $false, ". ", 2)

# ---------------------------------------------------------------------------
# 4. "input structure is constant over all samples" -> "input shape is consistent across all samples"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("input structure is constant over all samples", $true, $false, $false, $false, $false, $true, 1, $false, "input shape is consistent across all samples", 2)

# ---------------------------------------------------------------------------
# 5. Insert parenthetical about 3D CNNs after "videos" (before ", finally these features")
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("videos, finally these features are fed to", $true, $false, $false, $false, $false, $true, 1, $false, "videos (3D CNNs are particularly good here compared to CNN-LSTMs for lie detection as they can pick up on features related to motion in videos more which can allow the model to understand different gestures and facial movements that might be correlated with lying), finally these features are fed to", 2)

# ---------------------------------------------------------------------------
# 6. "data and discover more relations" -> "data and discover more hidden relations"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("data and discover more relations", $true, $false, $false, $false, $false, $true, 1, $false, "data and discover more hidden relations", 2)

# ---------------------------------------------------------------------------
# 7. Fix typo " 3DD" -> " 3D"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" 3DD CNN layers.", $true, $false, $false, $false, $false, $true, 1, $false, " 3D CNN layers.", 2)

# ---------------------------------------------------------------------------
# 8. Add sz/szCs 24 to the "Important remarks and discoveries" heading paragraph
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "Important remarks and discoveries") {
        $paraXml = '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Important remarks and discoveries</w:t></w:r></w:p>'
        $pkgXml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $paraXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $p.Range.InsertXML($pkgXml2)
    }
}

# ---------------------------------------------------------------------------
# 9 & 10. Rewrite the accuracy paragraph
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" reaching 95%", $true, $false, $false, $false, $false, $true, 1, $false, " reaching ~92%", 2)

$oldTail = " using a CNN-LSTM model, although it" + [char]0x2019 + "s a bit less with a 3D CNN at around 85%, but surprisingly when manually splitting the data to ensure that the test set contains videos of subjects that have never appeared in the train set (to eliminate all bias superstitions) the CNN-LSTM achieves poor results with accuracy below 60% while the 3D CNN maintains decent results with peek accuracy of 83% on both train and test sets."
$newTail = " using a CNN-LSTM model, although it" + [char]0x2019 + "s a lower with a 3D CNN at around 83%, but surprisingly when manually splitting the data to ensure that the test set contains videos of subjects that have never appeared in the train set (to eliminate all bias superstitions) the CNN-LSTM achieves poor results with accuracy below 60% while the 3D CNN maintains decent results with peek accuracy of around 83% on both train and test sets (same as the automatic, biased split) So we chose to continue with the 3D CNN as the preferred model when classifying the videos for visual cues that indicate lying."

$d.Content.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 2)
